# This script applies numeric corrections (refreshed market-board pricing
# and profit calculations) to several rows across the profession sheets of
# the Bahamut_Profits workbook, as produced by the scheduled pricing runner.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 2788.5293
$ws.Range("I113").Value = 3167.2222
$ws.Range("J113").Value = 2362.5
$ws.Range("K113").Value = 3167.2222
$ws.Range("L113").Value = 2362.5
$ws.Range("M113").Value = 86.77779999999984
$ws.Range("N113").Value = -8870.5

# Row 136
$ws.Range("H136").Value = 39533.332
$ws.Range("J136").Value = 39533.332
$ws.Range("L136").Value = 39533.332
$ws.Range("N136").Value = -49733.332

# Row 139
$ws.Range("H139").Value = 34060
$ws.Range("J139").Value = 34060
$ws.Range("L139").Value = 34060
$ws.Range("N139").Value = -44340

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 146.66667
$ws.Range("I4").Value = 146.66667
$ws.Range("K4").Value = 146.66667
$ws.Range("M4").Value = -30.66667000000001

# Row 5
$ws.Range("H5").Value = 72.5
$ws.Range("I5").Value = 45
$ws.Range("K5").Value = 45
$ws.Range("M5").Value = 67

# Row 61
$ws.Range("H61").Value = 4502.8
$ws.Range("I61").Value = 3500
$ws.Range("J61").Value = 5171.3335
$ws.Range("K61").Value = 3500
$ws.Range("L61").Value = 5171.3335
$ws.Range("M61").Value = -3288
$ws.Range("N61").Value = -5595.3335

# Row 63
$ws.Range("H63").Value = 1817.5
$ws.Range("I63").Value = 1781
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 1781
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1095
$ws.Range("N63").Value = -3372

# Row 66
$ws.Range("H66").Value = 1817.5
$ws.Range("I66").Value = 1781
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 8905
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -5473
$ws.Range("N66").Value = -16864

# Row 74
$ws.Range("H74").Value = 747.60974
$ws.Range("I74").Value = 787.2632
$ws.Range("J74").Value = 713.36365
$ws.Range("K74").Value = 787.2632
$ws.Range("L74").Value = 713.36365
$ws.Range("M74").Value = 86.73680000000002
$ws.Range("N74").Value = -2461.36365

# Row 77
$ws.Range("H77").Value = 747.60974
$ws.Range("I77").Value = 787.2632
$ws.Range("J77").Value = 713.36365
$ws.Range("K77").Value = 3936.316
$ws.Range("L77").Value = 3566.81825
$ws.Range("M77").Value = 431.6840000000002
$ws.Range("N77").Value = -12302.81825

# Row 132
$ws.Range("H132").Value = 4225.467
$ws.Range("I132").Value = 3737.3333
$ws.Range("K132").Value = 11211.9999
$ws.Range("M132").Value = -8681.999899999999

# Row 136
$ws.Range("H136").Value = 4502.8
$ws.Range("I136").Value = 3500
$ws.Range("J136").Value = 5171.3335
$ws.Range("K136").Value = 10500
$ws.Range("L136").Value = 15514.0005
$ws.Range("M136").Value = -7950
$ws.Range("N136").Value = -20614.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 72.5
$ws.Range("I4").Value = 45
$ws.Range("K4").Value = 45
$ws.Range("M4").Value = 70

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 2104.8462
$ws.Range("I132").Value = 1169.3334
$ws.Range("J132").Value = 2600.1177
$ws.Range("K132").Value = 3508.0002
$ws.Range("L132").Value = 7800.353099999999
$ws.Range("M132").Value = -978.0001999999999
$ws.Range("N132").Value = -12860.3531

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 112
$ws.Range("H112").Value = 3354.3845
$ws.Range("I112").Value = 902.3333
$ws.Range("J112").Value = 4090
$ws.Range("K112").Value = 2706.9999
$ws.Range("L112").Value = 12270
$ws.Range("M112").Value = -1598.9999
$ws.Range("N112").Value = -14486

# Row 117
$ws.Range("H117").Value = 1641.5834
$ws.Range("I117").Value = 1077
$ws.Range("J117").Value = 1754.5
$ws.Range("K117").Value = 3231
$ws.Range("L117").Value = 5263.5
$ws.Range("M117").Value = 211
$ws.Range("N117").Value = -12147.5

# Row 129
$ws.Range("H129").Value = 1513.05
$ws.Range("J129").Value = 1694.4286
$ws.Range("L129").Value = 5083.2858
$ws.Range("N129").Value = -15083.2858

# Row 131
$ws.Range("H131").Value = 38066.32
$ws.Range("J131").Value = 2783.7827
$ws.Range("L131").Value = 8351.348100000001
$ws.Range("N131").Value = -18431.3481

# Row 137
$ws.Range("H137").Value = 42230.37
$ws.Range("I137").Value = 2314.6667
$ws.Range("J137").Value = 92125
$ws.Range("K137").Value = 6944.000100000001
$ws.Range("L137").Value = 276375
$ws.Range("M137").Value = -1844.000100000001
$ws.Range("N137").Value = -286575

# Row 139
$ws.Range("H139").Value = 1706.5
$ws.Range("I139").Value = 1389.1
$ws.Range("K139").Value = 4167.299999999999
$ws.Range("M139").Value = 972.7000000000007

# Row 141
$ws.Range("H141").Value = 6124.4443
$ws.Range("I141").Value = 1665
$ws.Range("J141").Value = 9692
$ws.Range("K141").Value = 4995
$ws.Range("L141").Value = 29076
$ws.Range("M141").Value = 185
$ws.Range("N141").Value = -39436

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 53
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").ClearContents()

# Row 55
$ws.Range("H55").Value = 4000
$ws.Range("J55").Value = 4000
$ws.Range("L55").Value = 4000
$ws.Range("N55").Value = -4654

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1591.2727
$ws.Range("I61").Value = 1438
$ws.Range("K61").Value = 1438
$ws.Range("M61").Value = -1236

# Row 113
$ws.Range("H113").Value = 1591.2727
$ws.Range("I113").Value = 1438
$ws.Range("K113").Value = 1438
$ws.Range("M113").Value = 732

# Row 132
$ws.Range("H132").Value = 2319.568
$ws.Range("I132").Value = 2095.394
$ws.Range("J132").Value = 2992.0908
$ws.Range("K132").Value = 6286.181999999999
$ws.Range("L132").Value = 8976.2724
$ws.Range("M132").Value = -3756.181999999999
$ws.Range("N132").Value = -14036.2724

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 2813
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 3084
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 3084
$ws.Range("M96").Value = -627
$ws.Range("N96").Value = -5830

# Row 136
$ws.Range("H136").Value = 1881.8387
$ws.Range("I136").Value = 1711.6
$ws.Range("J136").Value = 2191.3635
$ws.Range("K136").Value = 5134.799999999999
$ws.Range("L136").Value = 6574.0905
$ws.Range("M136").Value = -2584.799999999999
$ws.Range("N136").Value = -11674.0905
